$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2328358208955224
$ws.Range("C2").Value = 0.5074626865671642
$ws.Range("J2").Value = 0.008955223880597015
$ws.Range("P2").Value = 0.1701492537313433
$ws.Range("S2").Value = 0.08059701492537313
$ws.Range("C3").Value = 0.05
$ws.Range("J3").Value = 0.005555555555555556
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.1444444444444444
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("O4").Value = 0.01724137931034483
$ws.Range("P4").Value = 0.6896551724137931
$ws.Range("S4").Value = 0.2241379310344828
$ws.Range("B6").Value = 0.06818181818181818
$ws.Range("D6").Value = 0.00909090909090909
$ws.Range("F6").Value = 0.04090909090909091
$ws.Range("J6").Value = 0.2818181818181818
$ws.Range("O6").Value = 0.01818181818181818
$ws.Range("Q6").Value = 0.1727272727272727
$ws.Range("R6").Value = 0.09545454545454546
$ws.Range("S6").Value = 0.3136363636363637
$ws.Range("B7").Value = 0.1326530612244898
$ws.Range("D7").Value = 0.01530612244897959
$ws.Range("F7").Value = 0.02551020408163265
$ws.Range("J7").Value = 0.173469387755102
$ws.Range("O7").Value = 0.02551020408163265
$ws.Range("Q7").Value = 0.1785714285714286
$ws.Range("R7").Value = 0.09693877551020408
$ws.Range("S7").Value = 0.3520408163265306
$ws.Range("B8").Value = 0.07875894988066826
$ws.Range("D8").Value = 0.02625298329355609
$ws.Range("F8").Value = 0.05727923627684964
$ws.Range("J8").Value = 0.1336515513126492
$ws.Range("O8").Value = 0.01909307875894988
$ws.Range("Q8").Value = 0.2004773269689737
$ws.Range("R8").Value = 0.1384248210023866
$ws.Range("S8").Value = 0.3460620525059666
$ws.Range("B9").Value = 0.108433734939759
$ws.Range("D9").Value = 0.02409638554216868
$ws.Range("E9").Value = 0.006024096385542169
$ws.Range("F9").Value = 0.06626506024096386
$ws.Range("J9").Value = 0.1506024096385542
$ws.Range("O9").Value = 0.03012048192771084
$ws.Range("Q9").Value = 0.1987951807228916
$ws.Range("R9").Value = 0.1265060240963855
$ws.Range("S9").Value = 0.2891566265060241
$ws.Range("B10").Value = 0.1170212765957447
$ws.Range("D10").Value = 0.02695035460992908
$ws.Range("E10").Value = 0.0007092198581560284
$ws.Range("F10").Value = 0.07021276595744681
$ws.Range("J10").Value = 0.1404255319148936
$ws.Range("O10").Value = 0.01205673758865248
$ws.Range("Q10").Value = 0.2347517730496454
$ws.Range("R10").Value = 0.09574468085106383
$ws.Range("S10").Value = 0.3021276595744681
$ws.Range("G11").Value = 0.1093247588424437
$ws.Range("J11").Value = 0.1286173633440514
$ws.Range("K11").Value = 0.1768488745980707
$ws.Range("L11").Value = 0.5787781350482315
$ws.Range("S11").Value = 0.006430868167202572
$ws.Range("G12").Value = 0.7081081081081081
$ws.Range("J12").Value = 0.2162162162162162
$ws.Range("K12").Value = 0.01621621621621622
$ws.Range("L12").Value = 0.03243243243243243
$ws.Range("S12").Value = 0.02702702702702703
$ws.Range("G13").Value = 0.7391304347826086
$ws.Range("J13").Value = 0.1739130434782609
$ws.Range("S13").Value = 0.08695652173913043
$ws.Range("F15").Value = 0.01363636363636364
$ws.Range("H15").Value = 0.1272727272727273
$ws.Range("I15").Value = 0.08636363636363636
$ws.Range("J15").Value = 0.3772727272727273
$ws.Range("K15").Value = 0.06818181818181818
$ws.Range("M15").Value = 0.01363636363636364
$ws.Range("O15").Value = 0.05454545454545454
$ws.Range("S15").Value = 0.2590909090909091
$ws.Range("F16").Value = 0.008547008547008548
$ws.Range("H16").Value = 0.1965811965811966
$ws.Range("I16").Value = 0.05982905982905983
$ws.Range("J16").Value = 0.4316239316239316
$ws.Range("K16").Value = 0.1025641025641026
$ws.Range("M16").Value = 0.02136752136752137
$ws.Range("O16").Value = 0.02136752136752137
$ws.Range("S16").Value = 0.1581196581196581
$ws.Range("F17").Value = 0.01346153846153846
$ws.Range("H17").Value = 0.1711538461538462
$ws.Range("I17").Value = 0.08461538461538462
$ws.Range("J17").Value = 0.4211538461538462
$ws.Range("K17").Value = 0.1019230769230769
$ws.Range("M17").Value = 0.02692307692307692
$ws.Range("O17").Value = 0.06923076923076923
$ws.Range("S17").Value = 0.1115384615384615
$ws.Range("F18").Value = 0.02766798418972332
$ws.Range("H18").Value = 0.1739130434782609
$ws.Range("I18").Value = 0.06719367588932806
$ws.Range("J18").Value = 0.4466403162055336
$ws.Range("K18").Value = 0.1383399209486166
$ws.Range("M18").Value = 0.01185770750988142
$ws.Range("O18").Value = 0.04743083003952569
$ws.Range("S18").Value = 0.08695652173913043
$ws.Range("F19").Value = 0.01790510295434199
$ws.Range("H19").Value = 0.1942703670546106
$ws.Range("J19").Value = 0.3849597135183527
$ws.Range("K19").Value = 0.1136974037600716
$ws.Range("M19").Value = 0.01969561324977619
$ws.Range("O19").Value = 0.08415398388540735
$ws.Range("S19").Value = 0.1208594449418084
